# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" scraped UI
# text block from the Requisitos section, together with the blank
# spacer paragraph and page-break paragraph that were generated along
# with it, collapsing the trailing paragraphs back down to the single
# page-break paragraph that should end the document.

$d = $word.ActiveDocument

$marker = "Ver no Jupiter Salvar em pdf Salvar em docx"

$searchRange = $d.Content
$found = $searchRange.Find.Execute($marker, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

if ($found) {
    # Paragraph that holds the matched text.
    $markerPara = $searchRange.Paragraphs.Item(1)
    $startIndex = $markerPara.Index

    # The marker paragraph is immediately followed by an empty
    # paragraph and then an empty page-break paragraph - all three
    # (plus the marker paragraph itself) need to go, leaving the
    # following paragraph (the real trailing page break) untouched.
    $firstParaToDelete = $d.Paragraphs.Item($startIndex)
    $lastParaToDelete = $d.Paragraphs.Item($startIndex + 3)

    $deleteRange = $d.Range($firstParaToDelete.Range.Start, $lastParaToDelete.Range.End)
    $deleteRange.Delete()
}
